# Insert a new weekly price record as row 305 in the "Vega Modelo de Temuco - Apio" data.
# This shifts the existing rows 305-323 down to 306-324 and fills the new row 305
# with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 305, pushing everything below it down by one.
$ws.Rows.Item(305).Insert()

# Populate the newly inserted row 305 with the new record's values.
$ws.Cells.Item(305, 1).Value2 = 10
$ws.Cells.Item(305, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(305, 3).Value2 = "La Araucanía"
$ws.Cells.Item(305, 4).Value2 = 44753
$ws.Cells.Item(305, 5).Value2 = 9
$ws.Cells.Item(305, 6).Value2 = 100112017
$ws.Cells.Item(305, 7).Value2 = "Apio"
$ws.Cells.Item(305, 8).Value2 = "Americana (o)"
$ws.Cells.Item(305, 9).Value2 = "Primera"
$ws.Cells.Item(305, 10).Value2 = 140
$ws.Cells.Item(305, 11).Value2 = 9000
$ws.Cells.Item(305, 12).Value2 = 9000
$ws.Cells.Item(305, 13).Value2 = 9000
$ws.Cells.Item(305, 14).Value2 = "`$/docena de matas"
$ws.Cells.Item(305, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(305, 16).Value2 = 1500
$ws.Cells.Item(305, 17).Value2 = 6
$ws.Cells.Item(305, 18).Value2 = "Hortaliza"

# Ensure the date cell keeps the same date/time number format as the other rows in column D.
$ws.Cells.Item(305, 4).NumberFormat = $ws.Cells.Item(306, 4).NumberFormat
